$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Role"
$ws.Range("C1").Value = "IDAM Roles"
$ws.Range("C1").Select()
